# Daily attendance processing - 2025-10-29 06:57:18
# Applies the day's attendance updates to the Session Analysis Results sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Narrow the "Status" column (I) a bit ---
$ws.Columns.Item(9).ColumnWidth = 9.17

# --- Row 2 (ANATOMY session 1): reorder the "Recorded By" list ---
$ws.Range("G2").Value = "Veronia.rafat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, gehanadel@med.asu.edu.eg, System"

# --- Summary box (K6:L10): Recorded/Missing sessions + coverage numbers ---
$ws.Range("L6").Value = 3
$ws.Range("L7").Value = 0

# --- Row 9 (HISTOLOGY session 1): reorder the "Recorded By" list ---
$ws.Range("G9").Value = "Shimaa.ashraf@med.asu.edu.eg, Safa.hany@med.asu.edu.eg"

# Coverage % / Average Attendance % - keep these as plain text (matching the
# original "NN.N%" formatting) rather than letting Excel auto-convert them
# into numeric percentages. We build the text via a literal-string formula
# and then flatten it back down to a plain value so no formula (and no
# extra number-format style) is left behind on the cell.
function Set-TextValue($addr, $text) {
  $r = $ws.Range($addr)
  $r.Formula = "=""" + $text + """"
  $r.Copy()
  $r.PasteSpecial(-4163)
}

Set-TextValue "L9" "10.3%"
Set-TextValue "L10" "10.2%"

# --- Row 15 / summary table (M15:S15) mirrors the K6:L10 box numbers ---
$ws.Range("O15").Value = 3
$ws.Range("P15").Value = 0
Set-TextValue "R15" "10.3%"
Set-TextValue "S15" "10.2%"
$excel.CutCopyMode = $false

# --- Row 28 (PHYSIOLOGY session 1): now recorded ---
# Pick up the green "Recorded" formatting (font/fill/alignment) used by the
# other recorded rows, e.g. row 9, by copying its formats across.
$ws.Range("A9:I9").Copy()
$ws.Range("A28:I28").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("G28").Value = "Aya_hamed@med.asu.edu.eg"
$ws.Range("H28").Value = "6/251"
$ws.Range("I28").Value = "Recorded"
